# Auto-generated edit script: updates crypto price/volume table cells
# to match the latest scrape (commit: "Updated cryptos list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, $CellRef, $NewValue)
    $cell = $Sheet.Range($CellRef)
    # Force the cell to Text format first so numeric-looking strings
    # (e.g. "212.51") are stored as text, not converted to numbers,
    # then clear the format again so no stray style is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.ClearFormats()
}

Set-TextCell $ws "D2" "29.605.82"
Set-TextCell $ws "E2" "  +3.43%  "

Set-TextCell $ws "D3" "1.607.45"
Set-TextCell $ws "E3" "  +2.78%  "

Set-TextCell $ws "D4" "1.00"
Set-TextCell $ws "E4" "  +0.02%  "

Set-TextCell $ws "D5" "212.51"

Set-TextCell $ws "D6" "0.520"
Set-TextCell $ws "E6" "  +2.79%  "

Set-TextCell $ws "E7" "  -0.01%  "

Set-TextCell $ws "D8" "26.88"
Set-TextCell $ws "E8" "  +7.98%  "

Set-TextCell $ws "D9" "43.55"
Set-TextCell $ws "E9" "  -1.43%  "

Set-TextCell $ws "E10" "  +2.64%  "

Set-TextCell $ws "E11" "  +2.56%  "

Set-TextCell $ws "E12" "  +1.42%  "

Set-TextCell $ws "D13" "1.836.43"
Set-TextCell $ws "E13" "  +2.73%  "

Set-TextCell $ws "D14" "1.615.89"
Set-TextCell $ws "E14" "  +3.37%  "

Set-TextCell $ws "D15" "29.613.60"
Set-TextCell $ws "E15" "  +3.33%  "

Set-TextCell $ws "E16" "  +3.94%  "

Set-TextCell $ws "E17" "  +2.41%  "

Set-TextCell $ws "D18" "63.40"
Set-TextCell $ws "E18" "  +3.15%  "

Set-TextCell $ws "D19" "240.88"
Set-TextCell $ws "E19" "  +6.15%  "

Set-TextCell $ws "D20" "7.59"
Set-TextCell $ws "E20" "  +3.93%  "

Set-TextCell $ws "E21" "  +1.91%  "

Set-TextCell $ws "E22" "  -0.05%  "

Set-TextCell $ws "E23" "  +1.87%  "

Set-TextCell $ws "D24" "9.23"
Set-TextCell $ws "E24" "  +2.34%  "

Set-TextCell $ws "E25" "  +0.63%  "

Set-TextCell $ws "D26" "154.53"
Set-TextCell $ws "E26" "  +1.98%  "

Set-TextCell $ws "B27" "EthereumClassic"
Set-TextCell $ws "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws "D27" "15.30"
Set-TextCell $ws "E27" "  +3.65%  "

Set-TextCell $ws "B28" "Stellar"
Set-TextCell $ws "C28" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws "D28" "0.108"
Set-TextCell $ws "E28" "  +3.27%  "

Set-TextCell $ws "D29" "6.41"
Set-TextCell $ws "E29" "  +3.26%  "

Set-TextCell $ws "E30" "  -0.04%  "

Set-TextCell $ws "E31" "  +3.46%  "

Set-TextCell $ws "E32" "  +0.89%  "

Set-TextCell $ws "E33" "  +1.59%  "

Set-TextCell $ws "E34" "  +4.34%  "

Set-TextCell $ws "D35" "1.408.66"
Set-TextCell $ws "E35" "  +0.64%  "

Set-TextCell $ws "E36" "  -0.05%  "

Set-TextCell $ws "E37" "  +5.07%  "

Set-TextCell $ws "E38" "  +5.32%  "

Set-TextCell $ws "E39" "  +0.14%  "

Set-TextCell $ws "D41" "0.539"
Set-TextCell $ws "E41" "  +4.55%  "

Set-TextCell $ws "E42" "  +2.55%  "

Set-TextCell $ws "D43" "0.0491"
Set-TextCell $ws "E43" "  +6.70%  "

Set-TextCell $ws "D44" "54.02"
Set-TextCell $ws "E44" "  +27.79%  "

Set-TextCell $ws "D45" "0.798"
Set-TextCell $ws "E45" "  +4.19%  "

Set-TextCell $ws "E46" "  -0.05%  "

Set-TextCell $ws "D47" "66.04"
Set-TextCell $ws "E47" "  +3.49%  "

Set-TextCell $ws "D48" "5.29"
Set-TextCell $ws "E48" "  +1.30%  "

Set-TextCell $ws "D49" "1.746.85"
Set-TextCell $ws "E49" "  +2.93%  "

Set-TextCell $ws "D50" "0.866"
Set-TextCell $ws "E50" "  -0.19%  "

Set-TextCell $ws "D51" "86.60"
Set-TextCell $ws "E51" "  +2.11%  "
